# Apply the "Update playlist and detail" change:
#  - Renumber the track-index column (A) on the "songs" sheet so that it
#    continues counting across albums (rows 8-40) instead of restarting
#    at 1 for every album.
#  - Update the stored selection on each sheet.

$wb = $excel.ActiveWorkbook

$wsAlbums = $wb.Worksheets.Item("albums")
$wsSongs  = $wb.Worksheets.Item("songs")

# Renumber column A for rows 8..40 on the songs sheet to a running count
# that continues on from row 7 (value 6).
$counter = 6
for ($row = 8; $row -le 40; $row++) {
    $counter = $counter + 1
    $wsSongs.Cells.Item($row, 1).Value = $counter
}

# Update the remembered selection on each sheet.
$wsAlbums.Range("B8").Select()
$wsSongs.Range("B3").Select()
